$wb = $excel.ActiveWorkbook

# ALC!8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 723.6
$ws.Range("I8").Value = 723.6
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2170.8
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2031.8
$ws.Range("N8").ClearContents()

# ALC!12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 254.5
$ws.Range("I12").Value = 299
$ws.Range("J12").Value = 210
$ws.Range("K12").Value = 299
$ws.Range("L12").Value = 210
$ws.Range("M12").Value = -129
$ws.Range("N12").Value = -550

# ALC!15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 481.25
$ws.Range("I15").Value = 481.25
$ws.Range("K15").Value = 1443.75
$ws.Range("M15").Value = -1274.75

# ALC!40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2285
$ws.Range("I40").Value = 1978.6
$ws.Range("J40").Value = 2591.4
$ws.Range("K40").Value = 1978.6
$ws.Range("L40").Value = 2591.4
$ws.Range("M40").Value = -1803.6
$ws.Range("N40").Value = -2941.4

# ALC!98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1018.9474
$ws.Range("I98").Value = 807.28
$ws.Range("J98").Value = 1426
$ws.Range("K98").Value = 807.28
$ws.Range("L98").Value = 1426
$ws.Range("M98").Value = 690.72
$ws.Range("N98").Value = -4422

# ALC!122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1018.9474
$ws.Range("I122").Value = 807.28
$ws.Range("J122").Value = 1426
$ws.Range("K122").Value = 2421.84
$ws.Range("L122").Value = 4278
$ws.Range("M122").Value = 28.15999999999985
$ws.Range("N122").Value = -9178

# ALC!129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 882.1429000000001
$ws.Range("J129").Value = 899.35486
$ws.Range("L129").Value = 2698.06458
$ws.Range("N129").Value = -12698.06458

# ALC!137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1791.5238
$ws.Range("I137").Value = 1479
$ws.Range("K137").Value = 4437
$ws.Range("M137").Value = -1887

# ALC!138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1912.5522
$ws.Range("I138").Value = 1726.2646
$ws.Range("J138").Value = 2104.4849
$ws.Range("K138").Value = 5178.793799999999
$ws.Range("L138").Value = 6313.4547
$ws.Range("M138").Value = -38.79379999999946
$ws.Range("N138").Value = -16593.4547

# ARM!74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1823.2667
$ws.Range("I74").Value = 832
$ws.Range("K74").Value = 832
$ws.Range("M74").Value = 42

# ARM!77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1823.2667
$ws.Range("I77").Value = 832
$ws.Range("K77").Value = 4160
$ws.Range("M77").Value = 208

# ARM!132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1428.6316
$ws.Range("I132").Value = 1428.6316
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4285.8948
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1755.8948
$ws.Range("N132").ClearContents()

# BSM!53
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

# BSM!94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1062.6364
$ws.Range("I94").Value = 931.8570999999999
$ws.Range("J94").Value = 1291.5
$ws.Range("K94").Value = 931.8570999999999
$ws.Range("L94").Value = 1291.5
$ws.Range("M94").Value = -480.8570999999999
$ws.Range("N94").Value = -2193.5

# BSM!99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1188.1428
$ws.Range("J99").Value = 1188.1428
$ws.Range("L99").Value = 1188.1428
$ws.Range("N99").Value = -4184.1428

# CRP!31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2515.25
$ws.Range("I31").Value = 3081.1667
$ws.Range("J31").Value = 1949.3334
$ws.Range("K31").Value = 3081.1667
$ws.Range("L31").Value = 1949.3334
$ws.Range("M31").Value = -2786.1667
$ws.Range("N31").Value = -2539.3334

# CRP!34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2515.25
$ws.Range("I34").Value = 3081.1667
$ws.Range("J34").Value = 1949.3334
$ws.Range("K34").Value = 3081.1667
$ws.Range("L34").Value = 1949.3334
$ws.Range("M34").Value = -2879.1667
$ws.Range("N34").Value = -2353.3334

# CRP!107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 723.3333
$ws.Range("J107").Value = 1385.25
$ws.Range("L107").Value = 1385.25
$ws.Range("N107").Value = -5225.25

# CRP!132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1847.9
$ws.Range("I132").Value = 1455.2858
$ws.Range("J132").Value = 2764
$ws.Range("K132").Value = 4365.857400000001
$ws.Range("L132").Value = 8292
$ws.Range("M132").Value = -1835.857400000001
$ws.Range("N132").Value = -13352

# CRP!134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2574.4443
$ws.Range("I134").Value = 2162.7334
$ws.Range("J134").Value = 4633
$ws.Range("K134").Value = 6488.2002
$ws.Range("L134").Value = 13899
$ws.Range("M134").Value = -3953.2002
$ws.Range("N134").Value = -18969

# CRP!135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# CRP!140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# CUL!2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 221.7
$ws.Range("I2").Value = 168
$ws.Range("J2").Value = 347
$ws.Range("K2").Value = 1008
$ws.Range("L2").Value = 2082
$ws.Range("M2").Value = -895
$ws.Range("N2").Value = -2308

# CUL!11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 933.3333
$ws.Range("I11").Value = 655.5
$ws.Range("J11").Value = 1489
$ws.Range("K11").Value = 1966.5
$ws.Range("L11").Value = 4467
$ws.Range("M11").Value = -1826.5
$ws.Range("N11").Value = -4747

# CUL!131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14241.528
$ws.Range("J131").Value = 14770.02
$ws.Range("L131").Value = 44310.06
$ws.Range("N131").Value = -54390.06

# GSM!80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -3996

# GSM!83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -19984

# LTW!134
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 49228
$ws.Range("J134").Value = 49228
$ws.Range("L134").Value = 49228
$ws.Range("N134").Value = -59368

# LTW!136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4339.5557
$ws.Range("I136").Value = 3517.5
$ws.Range("K136").Value = 10552.5
$ws.Range("M136").Value = -8002.5

# WVR!52
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 35000
$ws.Range("J52").Value = 35000
$ws.Range("L52").Value = 35000
$ws.Range("N52").Value = -35452

# WVR!100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1076.8334
$ws.Range("I100").Value = 932.2
$ws.Range("K100").Value = 1864.4
$ws.Range("M100").Value = -1323.4

# WVR!136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 17923244
$ws.Range("I136").Value = 30866366
$ws.Range("J136").Value = 1999.4615
$ws.Range("K136").Value = 92599098
$ws.Range("L136").Value = 5998.3845
$ws.Range("M136").Value = -92596548
$ws.Range("N136").Value = -11098.3845
